$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1 (Código de barras (EAN)) and G1 (NCM) ---
# Copy the existing header formatting (bold font, border, centered/top
# alignment) from E1 so the two new headers match the rest of the row.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "Código de barras (EAN)"
$ws.Range("G1").Value = "NCM"

# --- Row 2 (C-5682) new data ---
# EAN / NCM codes are entered as text (leading apostrophe forces text
# entry for the purely-numeric looking values), then the quote-prefix
# cell style that Excel applies is cleared back to Normal so the cells
# keep plain "General" formatting like the rest of the data rows.
$ws.Range("F2").Value = "'7893049568296"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'87083090"
$ws.Range("G2").Style = "Normal"

# --- Row 3 (HG 41123) new data ---
$ws.Range("F3").Value = "'7890903099835"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "8708.80.00"

# New column width for the NCM column
$ws.Columns.Item(7).ColumnWidth = 16.71

# Move the active cell/selection to D7
[void]$ws.Range("D7").Select()
